$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The window had been saved minimized; restore it to a normal state.
$wb.Windows.Item(1).WindowState = -4143  # xlNormal

# Add a new "Blood Bank" login row under the existing table (row 16).
$ws.Range("C16").Value = "BB_user"
$ws.Range("D16").Value = "BB_user12345"
$ws.Range("E16").Value = "Blood Bank"

# Give the new row the same thin left/right cell-dividers as the rest of
# the table (build the combined border on the first cell, then fan it
# out to the remaining cells in the row via a format copy).
$first = $ws.Range("C16")
$first.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$first.Borders.Item(7).Weight = 2      # xlThin
$first.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$first.Borders.Item(10).Weight = 2     # xlThin

$first.Copy()
$ws.Range("D16:E16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Leave the same cell selected that was active after the edit.
$ws.Range("E16").Select()
